$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Mark "Remove silly option thing (sidemenu)" (row 9) as Done.
#    The side menu / options menu is being disabled because it is
#    non-functional (it disappears as soon as the mouse moves).
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "Done"

# ---------------------------------------------------------------------------
# 2. Insert two new blank rows just above the final "Other menus" section
#    (old row 26), pushing it down to row 28, and leaving rows 26/27 blank.
# ---------------------------------------------------------------------------
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).Insert()

# ---------------------------------------------------------------------------
# 3. Populate the new row 25 with the new TODO entry about the Guide
#    control scrolling direction.
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "Guide should roll up and down with the mouse wheel, not left to right. "
$ws.Range("C25").Value = "Can do in a keymap, but not sure how in a skin. Requested."

# Match the look of the other two-line (wrapped) rows such as row 6/10/13.
$ws.Range("A25").WrapText = $true
$ws.Range("C25").WrapText = $true
$ws.Range("A25").VerticalAlignment = -4160
$ws.Range("C25").VerticalAlignment = -4160
$ws.Rows.Item(25).RowHeight = 28.8

# ---------------------------------------------------------------------------
# 4. Align the "Status" column (B) cells to the top of the row -- these rows
#    wrap to two lines so the short "Done" labels look better pinned to the
#    top rather than vertically centred.
# ---------------------------------------------------------------------------
$statusCells = @("B1", "B2", "B3", "B4", "B6", "B7", "B8", "B9", "B13", "B14", "B15", "B17", "B20")
foreach ($cellRef in $statusCells) {
    $ws.Range($cellRef).VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 5. Update the view: scroll so row 14 is at the top and the last edited
#    cell (C26) is selected, matching where work on the sheet left off.
# ---------------------------------------------------------------------------
$ws.Range("C26").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
